# Generate Report for Handback
# ------------------------------------------------------------------
# 1) "Ready for handoff" -> "Handed back: in sync with en-US"
#    (this shared string is used by Overview!B2/C2/B3/C3 and by
#    zh-cn!C2/C3 and de-de!C2/C3 - update every occurrence so the
#    shared text changes everywhere, matching a shared-string edit)
# 2) zh-cn!H2 / zh-cn!H3 ("Latest Handback DateTime"): the placeholder
#    "0001-01-01 00:00:00" becomes the real handback time
#    "2016-03-11 16:43:27"
# 3) de-de!H2 / de-de!H3 get a *different*, later handback time
#    "2016-03-11 16:43:33"
# 4) Populate the new "Latest Target File" (F) and "Latest Handback
#    File" (G) columns on both the zh-cn and de-de sheets, with the
#    same display text + hyperlink target as the existing "Source
#    File Name" (A2) / "Latest Handoff File" (D2) hyperlinks.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# --- Latest Handback DateTime -------------------------------------
$wsZh.Range("H2").Value = "2016-03-11 16:43:27"
$wsZh.Range("H3").Value = "2016-03-11 16:43:27"

$wsDe.Range("H2").Value = "2016-03-11 16:43:33"
$wsDe.Range("H3").Value = "2016-03-11 16:43:33"

# --- Latest Target File (F) / Latest Handback File (G) ------------
$mdDisplayZh = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.md"
$mdTargetZh  = "https://github.com/OpenLocalizationTest/oltest/blob/f5f7084dceee8714fda4d34f3740530fd7eef2c5/e2e/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.md"

$xlfDisplayZh = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.zh-cn.xlf"
$xlfTargetZh  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dd0c65b88bfabe45b8eb04d2ef09e11d023c83bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.zh-cn.xlf"

$xlfDisplayDe = "e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.de-de.xlf"
$xlfTargetDe  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/015e05080a9555e922435f75d8b7f5fd85685a8a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e9c65714-9c9e-4b6b-b1bb-fcde389cc5b8.dc77bb0a579a94acbd08830e9bd31675608e4439.de-de.xlf"

function Add-ReportLink($ws, $cellRef, $target, $display) {
    $rng = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($rng, $target, "", "", $display) | Out-Null
    $rng.Style = "Hyperlink"
}

# zh-cn sheet: rows 2 and 3
Add-ReportLink $wsZh "F2" $mdTargetZh  $mdDisplayZh
Add-ReportLink $wsZh "G2" $xlfTargetZh $xlfDisplayZh
Add-ReportLink $wsZh "F3" $mdTargetZh  $mdDisplayZh
Add-ReportLink $wsZh "G3" $xlfTargetZh $xlfDisplayZh

# de-de sheet: rows 2 and 3
Add-ReportLink $wsDe "F2" $mdTargetZh  $mdDisplayZh
Add-ReportLink $wsDe "G2" $xlfTargetDe $xlfDisplayDe
Add-ReportLink $wsDe "F3" $mdTargetZh  $mdDisplayZh
Add-ReportLink $wsDe "G3" $xlfTargetDe $xlfDisplayDe
